# Update incidence prediction values (AVG/L95CI/U95CI for Retinopatia de fondo INC)
# resulting from retraining the RandomForestRegressor model.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 70.20704477611932
$ws.Range("B2").Value = 51.64800000000002
$ws.Range("C2").Value = 82.76000000000001
$ws.Range("A3").Value = 0.8654925373134325
$ws.Range("B3").Value = 0.07200000000000002
$ws.Range("C3").Value = 1.879999999999999
$ws.Range("A4").Value = 3.7212935323383
$ws.Range("B4").Value = 1.38
$ws.Range("C4").Value = 6.496
$ws.Range("A5").Value = 0.9647761194029849
$ws.Range("B5").Value = 0.052
$ws.Range("C5").Value = 2.108000000000001
$ws.Range("A6").Value = 16.90185074626865
$ws.Range("B6").Value = 8.375999999999998
$ws.Range("C6").Value = 25.04399999999999
$ws.Range("A7").Value = 66.60573134328347
$ws.Range("B7").Value = 49.256
$ws.Range("C7").Value = 80.03199999999997
$ws.Range("A8").Value = 96.39209950248747
$ws.Range("B8").Value = 93.49600000000005
$ws.Range("C8").Value = 98.428
$ws.Range("A9").Value = 15.28031840796019
$ws.Range("B9").Value = 7.595999999999998
$ws.Range("C9").Value = 23.66399999999999
$ws.Range("A10").Value = 96.36539303482574
$ws.Range("B10").Value = 93.40800000000007
$ws.Range("C10").Value = 98.22799999999992
$ws.Range("A11").Value = 86.3321791044775
$ws.Range("B11").Value = 75.33200000000005
$ws.Range("C11").Value = 93.10399999999996
$ws.Range("A12").Value = 45.08495522388047
$ws.Range("B12").Value = 27.92000000000001
$ws.Range("C12").Value = 59.71199999999999
$ws.Range("A13").Value = 95.46189054726355
$ws.Range("B13").Value = 91.56400000000004
$ws.Range("C13").Value = 97.94400000000002
$ws.Range("A14").Value = 92.53745273631829
$ws.Range("B14").Value = 86.18400000000001
$ws.Range("C14").Value = 96.63600000000004
$ws.Range("A15").Value = 1.163044776119401
$ws.Range("B15").Value = 0.136
$ws.Range("C15").Value = 2.36
$ws.Range("A16").Value = 17.00549253731343
$ws.Range("B16").Value = 8.535999999999998
$ws.Range("C16").Value = 25.07599999999999
$ws.Range("A17").Value = 44.61460696517405
$ws.Range("B17").Value = 27.45199999999999
$ws.Range("C17").Value = 59.58799999999996
$ws.Range("A18").Value = 3.7212935323383
$ws.Range("B18").Value = 1.38
$ws.Range("C18").Value = 6.496
$ws.Range("A19").Value = 86.56181094527358
$ws.Range("B19").Value = 77.14799999999997
$ws.Range("C19").Value = 93.44799999999998
$ws.Range("A20").Value = 18.68131343283582
$ws.Range("B20").Value = 9.415999999999997
$ws.Range("C20").Value = 27.34
$ws.Range("A21").Value = 5.135641791044767
$ws.Range("B21").Value = 2.352000000000002
$ws.Range("C21").Value = 8.620000000000006
$ws.Range("A22").Value = 96.48328358208944
$ws.Range("B22").Value = 94.02399999999996
$ws.Range("C22").Value = 98.37999999999998
$ws.Range("A23").Value = 1.196218905472636
$ws.Range("B23").Value = 0.1999999999999999
$ws.Range("C23").Value = 2.615999999999998
$ws.Range("A24").Value = 16.9916815920398
$ws.Range("B24").Value = 8.555999999999997
$ws.Range("C24").Value = 25.21199999999999
$ws.Range("A25").Value = 14.27719402985073
$ws.Range("B25").Value = 7.275999999999995
$ws.Range("C25").Value = 22.05600000000001
$ws.Range("A26").Value = 5.094228855721386
$ws.Range("B26").Value = 2.204000000000002
$ws.Range("C26").Value = 8.672000000000008
$ws.Range("A27").Value = 95.51659701492525
$ws.Range("B27").Value = 92.00000000000001
$ws.Range("C27").Value = 98.21199999999995
$ws.Range("A28").Value = 68.60244776119393
$ws.Range("B28").Value = 50.52800000000001
$ws.Range("C28").Value = 81.73600000000006
$ws.Range("A29").Value = 1.228278606965172
$ws.Range("B29").Value = 0.2039999999999999
$ws.Range("C29").Value = 2.643999999999998
$ws.Range("A30").Value = 43.14923383084566
$ws.Range("B30").Value = 25.596
$ws.Range("C30").Value = 57.70399999999998
$ws.Range("A31").Value = 92.3635621890546
$ws.Range("B31").Value = 85.97599999999998
$ws.Range("C31").Value = 96.43600000000001
$ws.Range("A32").Value = 95.34121393034813
$ws.Range("B32").Value = 91.43600000000002
$ws.Range("C32").Value = 97.81600000000007
$ws.Range("A33").Value = 92.68531343283577
$ws.Range("B33").Value = 86.00800000000007
$ws.Range("C33").Value = 96.536
$ws.Range("A34").Value = 72.33343283582083
$ws.Range("B34").Value = 56.10799999999998
$ws.Range("C34").Value = 85.28399999999993
$ws.Range("A35").Value = 17.54883582089551
$ws.Range("B35").Value = 9.760000000000002
$ws.Range("C35").Value = 26.56399999999999
$ws.Range("A36").Value = 66.90963184079595
$ws.Range("B36").Value = 48.9
$ws.Range("C36").Value = 80.44000000000004
$ws.Range("A37").Value = 5.575363184079588
$ws.Range("B37").Value = 2.524000000000001
$ws.Range("C37").Value = 9.079999999999997
$ws.Range("A38").Value = 92.74728358208949
$ws.Range("B38").Value = 86.66000000000004
$ws.Range("C38").Value = 96.62
$ws.Range("A39").Value = 45.19542288557196
$ws.Range("B39").Value = 27.93599999999999
$ws.Range("C39").Value = 60.53199999999996
$ws.Range("A40").Value = 4.311402985074625
$ws.Range("B40").Value = 1.903999999999999
$ws.Range("C40").Value = 7.451999999999998
$ws.Range("A41").Value = 63.99906467661678
$ws.Range("B41").Value = 45.72000000000003
$ws.Range("C41").Value = 77.81599999999999
$ws.Range("A42").Value = 86.50762189054721
$ws.Range("B42").Value = 75.59600000000002
$ws.Range("C42").Value = 93.19999999999993
$ws.Range("A43").Value = 43.14923383084566
$ws.Range("B43").Value = 25.596
$ws.Range("C43").Value = 57.70399999999998
$ws.Range("A44").Value = 84.24352238805965
$ws.Range("B44").Value = 71.69200000000001
$ws.Range("C44").Value = 92.07999999999997
$ws.Range("A45").Value = 86.74009950248757
$ws.Range("B45").Value = 76.77599999999993
$ws.Range("C45").Value = 93.34
$ws.Range("A46").Value = 92.19359203980095
$ws.Range("B46").Value = 85.72800000000002
$ws.Range("C46").Value = 96.24799999999998
$ws.Range("A47").Value = 96.4461691542287
$ws.Range("B47").Value = 93.62800000000007
$ws.Range("C47").Value = 98.35199999999995
$ws.Range("A48").Value = 45.22314427860688
$ws.Range("B48").Value = 28.12
$ws.Range("C48").Value = 60.34399999999998
$ws.Range("A49").Value = 4.10666666666665
$ws.Range("B49").Value = 1.759999999999999
$ws.Range("C49").Value = 7.008000000000002
$ws.Range("A50").Value = 5.562348258706453
$ws.Range("B50").Value = 2.408000000000002
$ws.Range("C50").Value = 9.224000000000002
$ws.Range("A51").Value = 0.8154427860696513
$ws.Range("B51").Value = 0.112
$ws.Range("C51").Value = 1.871999999999999
$ws.Range("A52").Value = 92.45902487562179
$ws.Range("B52").Value = 85.81200000000001
$ws.Range("C52").Value = 96.30399999999992
$ws.Range("A53").Value = 45.03976119402976
$ws.Range("B53").Value = 27.99999999999999
$ws.Range("C53").Value = 60.65599999999995
$ws.Range("A54").Value = 1.148179104477611
$ws.Range("B54").Value = 0.1360000000000001
$ws.Range("C54").Value = 2.48
$ws.Range("A55").Value = 5.080557213930341
$ws.Range("B55").Value = 2.260000000000002
$ws.Range("C55").Value = 8.592000000000006
$ws.Range("A56").Value = 92.67313432835812
$ws.Range("B56").Value = 86.29200000000004
$ws.Range("C56").Value = 96.47199999999999
$ws.Range("A57").Value = 95.30473631840783
$ws.Range("B57").Value = 91.50400000000008
$ws.Range("C57").Value = 97.84000000000007
$ws.Range("A58").Value = 5.561014925373122
$ws.Range("B58").Value = 2.444000000000001
$ws.Range("C58").Value = 9.180000000000001
$ws.Range("A59").Value = 85.69221890547247
$ws.Range("B59").Value = 74.76799999999994
$ws.Range("C59").Value = 92.68000000000001
$ws.Range("A60").Value = 1.197134328358207
$ws.Range("B60").Value = 0.2280000000000001
$ws.Range("C60").Value = 2.479999999999999
$ws.Range("A61").Value = 92.74455721393034
$ws.Range("B61").Value = 85.73600000000006
$ws.Range("C61").Value = 96.58800000000005
$ws.Range("A62").Value = 86.30479601990041
$ws.Range("B62").Value = 74.96800000000005
$ws.Range("C62").Value = 93.07999999999994
$ws.Range("A63").Value = 92.78796019900484
$ws.Range("B63").Value = 86.72799999999998
$ws.Range("C63").Value = 96.68000000000004
$ws.Range("A64").Value = 1.197492537313431
$ws.Range("B64").Value = 0.1720000000000001
$ws.Range("C64").Value = 2.479999999999999
$ws.Range("A65").Value = 95.15235820895521
$ws.Range("B65").Value = 90.52799999999999
$ws.Range("C65").Value = 97.572
$ws.Range("A66").Value = 72.14300497512437
$ws.Range("B66").Value = 55.912
$ws.Range("C66").Value = 84.43199999999997
$ws.Range("A67").Value = 70.93996019900486
$ws.Range("B67").Value = 53.44400000000005
$ws.Range("C67").Value = 84.072
$ws.Range("A68").Value = 4.708159203980097
$ws.Range("B68").Value = 1.975999999999999
$ws.Range("C68").Value = 7.896000000000006
$ws.Range("A69").Value = 96.51301492537299
$ws.Range("B69").Value = 93.848
$ws.Range("C69").Value = 98.42399999999999
$ws.Range("A70").Value = 4.370348258706458
$ws.Range("B70").Value = 1.843999999999999
$ws.Range("C70").Value = 7.520000000000002
$ws.Range("A71").Value = 34.68694527363174
$ws.Range("B71").Value = 19.92400000000001
$ws.Range("C71").Value = 48.34000000000001
$ws.Range("A72").Value = 39.12135323383077
$ws.Range("B72").Value = 23.17599999999998
$ws.Range("C72").Value = 53.876
